$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17; existing rows 17..131 shift down to 18..132.
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with the new week's data.
$ws.Cells.Item(17, 1).Value = 4
$ws.Cells.Item(17, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(17, 3).Value = "Los Lagos"
$ws.Cells.Item(17, 4).Value = 44503
$ws.Cells.Item(17, 5).Value = 10
$ws.Cells.Item(17, 6).Value = 100112039
$ws.Cells.Item(17, 7).Value = "Ciboulette"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 40
$ws.Cells.Item(17, 11).Value = 2500
$ws.Cells.Item(17, 12).Value = 2500
$ws.Cells.Item(17, 13).Value = 2500
$ws.Cells.Item(17, 14).Value = "`$/docena de atados"
$ws.Cells.Item(17, 15).Value = "Región Metropolitana"
$ws.Cells.Item(17, 16).Value = 833
$ws.Cells.Item(17, 17).Value = 3
$ws.Cells.Item(17, 18).Value = "Hortaliza"
